$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.018117615953087807
$ws.Range("C2").Value = 0.008500768803060055
$ws.Range("D2").Value = 0.005139508750289679
$ws.Range("E2").Value = 0.003978007007390261
$ws.Range("F2").Value = 0.00025925389491021633
$ws.Range("I2").Value = 1.2575732469558716
$ws.Range("J2").Value = 0.1273263394832611
$ws.Range("K2").Value = 1.4221693277359009
